$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; this shifts existing rows 45-89 down to 46-90.
$ws.Rows.Item(45).Insert()

# Populate the new row 45 with the data.
$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(45, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(45, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(45, 4).Value = 44629
$ws.Cells.Item(45, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(45, 5).Value = 15
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100102
$ws.Cells.Item(45, 8).Value = "Cítricos"
$ws.Cells.Item(45, 9).Value = 100102004
$ws.Cells.Item(45, 10).Value = "Mandarina"
$ws.Cells.Item(45, 11).Value = "Murcott"
$ws.Cells.Item(45, 12).Value = "Tercera"
$ws.Cells.Item(45, 13).Value = 250
$ws.Cells.Item(45, 14).Value = 13000
$ws.Cells.Item(45, 15).Value = 14000
$ws.Cells.Item(45, 16).Value = 13500
$ws.Cells.Item(45, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(45, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(45, 19).Value = 675
$ws.Cells.Item(45, 20).Value = 20
